$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.467623114585876
$ws.Range("B1").Value = 3.811535358428955
$ws.Range("C1").Value = 5.672089576721191
$ws.Range("D1").Value = 1.638223648071289
$ws.Range("E1").Value = 0.9616954922676086
